# Auto-generated edit script for optimisation_result.xlsx (run 201 update)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Schedule
$ws2 = $wb.Worksheets.Item(2)   # Detailed

# ---- Sheet1 "Schedule": refresh totals on row 2, append row 3 ----
$ws1.Cells.Item(2, 5).Value = 582.3059580000001   # Cost ($)
$ws1.Cells.Item(2, 6).Value = 14.0044722943723   # Unit Cost ($/ML)

$ws1.Cells.Item(3, 1).Value = 46050.20833333334   # Start Time
$ws1.Cells.Item(3, 2).Value = 46050.66666666666   # Stop Time
$ws1.Cells.Item(3, 3).Value = 11.0   # Duration (h)
$ws1.Cells.Item(3, 4).Value = 41.58   # Volume (ML)
$ws1.Cells.Item(3, 5).Value = 443.1063877499999   # Cost ($)
$ws1.Cells.Item(3, 6).Value = 10.65671928210678   # Unit Cost ($/ML)
$ws1.Cells.Item(3, 1).NumberFormat = $ws1.Cells.Item(2, 1).NumberFormat
$ws1.Cells.Item(3, 2).NumberFormat = $ws1.Cells.Item(2, 2).NumberFormat

# ---- Sheet2 "Detailed": revise Price/Type for existing rows 13-49 ----
$ws2.Cells.Item(13, 2).Value = 78.00008
$ws2.Cells.Item(14, 2).Value = 84.79
$ws2.Cells.Item(15, 2).Value = 57.49018
$ws2.Cells.Item(15, 3).Value = "historical"
$ws2.Cells.Item(16, 2).Value = 56.98
$ws2.Cells.Item(16, 3).Value = "historical"
$ws2.Cells.Item(17, 3).Value = "historical"
$ws2.Cells.Item(18, 3).Value = "historical"
$ws2.Cells.Item(19, 3).Value = "historical"
$ws2.Cells.Item(20, 2).Value = 22.07
$ws2.Cells.Item(20, 3).Value = "historical"
$ws2.Cells.Item(21, 2).Value = 36.06
$ws2.Cells.Item(21, 3).Value = "historical"
$ws2.Cells.Item(22, 3).Value = "historical"
$ws2.Cells.Item(23, 2).Value = 22.84723
$ws2.Cells.Item(23, 3).Value = "historical"
$ws2.Cells.Item(24, 2).Value = 0.51
$ws2.Cells.Item(24, 3).Value = "historical"
$ws2.Cells.Item(25, 2).Value = 0.51
$ws2.Cells.Item(25, 3).Value = "historical"
$ws2.Cells.Item(26, 3).Value = "historical"
$ws2.Cells.Item(27, 3).Value = "historical"
$ws2.Cells.Item(28, 2).Value = 7.99033
$ws2.Cells.Item(28, 3).Value = "historical"
$ws2.Cells.Item(29, 3).Value = "historical"
$ws2.Cells.Item(30, 3).Value = "historical"
$ws2.Cells.Item(31, 2).Value = 36.06
$ws2.Cells.Item(31, 3).Value = "historical"
$ws2.Cells.Item(32, 2).Value = 0.7
$ws2.Cells.Item(32, 3).Value = "historical"
$ws2.Cells.Item(33, 2).Value = 0.51
$ws2.Cells.Item(34, 2).Value = 35.88
$ws2.Cells.Item(35, 2).Value = 36.06
$ws2.Cells.Item(36, 2).Value = 36.00466
$ws2.Cells.Item(37, 2).Value = 27.12548
$ws2.Cells.Item(38, 2).Value = 21.16499
$ws2.Cells.Item(39, 2).Value = 71.38292
$ws2.Cells.Item(40, 2).Value = 103.52221
$ws2.Cells.Item(41, 2).Value = 107.40667
$ws2.Cells.Item(42, 2).Value = 108.89
$ws2.Cells.Item(44, 2).Value = 98.4973
$ws2.Cells.Item(45, 2).Value = 87.93914
$ws2.Cells.Item(46, 2).Value = 89.00783
$ws2.Cells.Item(47, 2).Value = 101.25
$ws2.Cells.Item(48, 2).Value = 101.25
$ws2.Cells.Item(49, 2).Value = 95.91848

# ---- Sheet2 "Detailed": append new rows 50-97 (run 201 forecast horizon) ----
# Row 50
$ws2.Cells.Item(50, 1).Value = 46050.0
$ws2.Cells.Item(50, 2).Value = 84.79
$ws2.Cells.Item(50, 3).Value = "forecast"
$ws2.Cells.Item(50, 4).Value = 46050.0
$ws2.Cells.Item(50, 5).Value = "OFF"
$ws2.Cells.Item(50, 1).NumberFormat = $ws2.Cells.Item(49, 1).NumberFormat
$ws2.Cells.Item(50, 4).NumberFormat = $ws2.Cells.Item(49, 4).NumberFormat

# Row 51
$ws2.Cells.Item(51, 1).Value = 46050.02083333334
$ws2.Cells.Item(51, 2).Value = 92.45217
$ws2.Cells.Item(51, 3).Value = "forecast"
$ws2.Cells.Item(51, 4).Value = 46050.0
$ws2.Cells.Item(51, 5).Value = "OFF"
$ws2.Cells.Item(51, 1).NumberFormat = $ws2.Cells.Item(50, 1).NumberFormat
$ws2.Cells.Item(51, 4).NumberFormat = $ws2.Cells.Item(50, 4).NumberFormat

# Row 52
$ws2.Cells.Item(52, 1).Value = 46050.04166666666
$ws2.Cells.Item(52, 2).Value = 84.79
$ws2.Cells.Item(52, 3).Value = "forecast"
$ws2.Cells.Item(52, 4).Value = 46050.0
$ws2.Cells.Item(52, 5).Value = "OFF"
$ws2.Cells.Item(52, 1).NumberFormat = $ws2.Cells.Item(51, 1).NumberFormat
$ws2.Cells.Item(52, 4).NumberFormat = $ws2.Cells.Item(51, 4).NumberFormat

# Row 53
$ws2.Cells.Item(53, 1).Value = 46050.0625
$ws2.Cells.Item(53, 2).Value = 78.0
$ws2.Cells.Item(53, 3).Value = "forecast"
$ws2.Cells.Item(53, 4).Value = 46050.0
$ws2.Cells.Item(53, 5).Value = "OFF"
$ws2.Cells.Item(53, 1).NumberFormat = $ws2.Cells.Item(52, 1).NumberFormat
$ws2.Cells.Item(53, 4).NumberFormat = $ws2.Cells.Item(52, 4).NumberFormat

# Row 54
$ws2.Cells.Item(54, 1).Value = 46050.08333333334
$ws2.Cells.Item(54, 2).Value = 78.0
$ws2.Cells.Item(54, 3).Value = "forecast"
$ws2.Cells.Item(54, 4).Value = 46050.0
$ws2.Cells.Item(54, 5).Value = "OFF"
$ws2.Cells.Item(54, 1).NumberFormat = $ws2.Cells.Item(53, 1).NumberFormat
$ws2.Cells.Item(54, 4).NumberFormat = $ws2.Cells.Item(53, 4).NumberFormat

# Row 55
$ws2.Cells.Item(55, 1).Value = 46050.10416666666
$ws2.Cells.Item(55, 2).Value = 78.0
$ws2.Cells.Item(55, 3).Value = "forecast"
$ws2.Cells.Item(55, 4).Value = 46050.0
$ws2.Cells.Item(55, 5).Value = "OFF"
$ws2.Cells.Item(55, 1).NumberFormat = $ws2.Cells.Item(54, 1).NumberFormat
$ws2.Cells.Item(55, 4).NumberFormat = $ws2.Cells.Item(54, 4).NumberFormat

# Row 56
$ws2.Cells.Item(56, 1).Value = 46050.125
$ws2.Cells.Item(56, 2).Value = 78.00015
$ws2.Cells.Item(56, 3).Value = "forecast"
$ws2.Cells.Item(56, 4).Value = 46050.0
$ws2.Cells.Item(56, 5).Value = "OFF"
$ws2.Cells.Item(56, 1).NumberFormat = $ws2.Cells.Item(55, 1).NumberFormat
$ws2.Cells.Item(56, 4).NumberFormat = $ws2.Cells.Item(55, 4).NumberFormat

# Row 57
$ws2.Cells.Item(57, 1).Value = 46050.14583333334
$ws2.Cells.Item(57, 2).Value = 78.00015
$ws2.Cells.Item(57, 3).Value = "forecast"
$ws2.Cells.Item(57, 4).Value = 46050.0
$ws2.Cells.Item(57, 5).Value = "OFF"
$ws2.Cells.Item(57, 1).NumberFormat = $ws2.Cells.Item(56, 1).NumberFormat
$ws2.Cells.Item(57, 4).NumberFormat = $ws2.Cells.Item(56, 4).NumberFormat

# Row 58
$ws2.Cells.Item(58, 1).Value = 46050.16666666666
$ws2.Cells.Item(58, 2).Value = 79.82787
$ws2.Cells.Item(58, 3).Value = "forecast"
$ws2.Cells.Item(58, 4).Value = 46050.0
$ws2.Cells.Item(58, 5).Value = "OFF"
$ws2.Cells.Item(58, 1).NumberFormat = $ws2.Cells.Item(57, 1).NumberFormat
$ws2.Cells.Item(58, 4).NumberFormat = $ws2.Cells.Item(57, 4).NumberFormat

# Row 59
$ws2.Cells.Item(59, 1).Value = 46050.1875
$ws2.Cells.Item(59, 2).Value = 80.53772
$ws2.Cells.Item(59, 3).Value = "forecast"
$ws2.Cells.Item(59, 4).Value = 46050.0
$ws2.Cells.Item(59, 5).Value = "OFF"
$ws2.Cells.Item(59, 1).NumberFormat = $ws2.Cells.Item(58, 1).NumberFormat
$ws2.Cells.Item(59, 4).NumberFormat = $ws2.Cells.Item(58, 4).NumberFormat

# Row 60
$ws2.Cells.Item(60, 1).Value = 46050.20833333334
$ws2.Cells.Item(60, 2).Value = 84.79
$ws2.Cells.Item(60, 3).Value = "forecast"
$ws2.Cells.Item(60, 4).Value = 46050.0
$ws2.Cells.Item(60, 5).Value = "ON"
$ws2.Cells.Item(60, 1).NumberFormat = $ws2.Cells.Item(59, 1).NumberFormat
$ws2.Cells.Item(60, 4).NumberFormat = $ws2.Cells.Item(59, 4).NumberFormat

# Row 61
$ws2.Cells.Item(61, 1).Value = 46050.22916666666
$ws2.Cells.Item(61, 2).Value = 97.74777
$ws2.Cells.Item(61, 3).Value = "forecast"
$ws2.Cells.Item(61, 4).Value = 46050.0
$ws2.Cells.Item(61, 5).Value = "ON"
$ws2.Cells.Item(61, 1).NumberFormat = $ws2.Cells.Item(60, 1).NumberFormat
$ws2.Cells.Item(61, 4).NumberFormat = $ws2.Cells.Item(60, 4).NumberFormat

# Row 62
$ws2.Cells.Item(62, 1).Value = 46050.25
$ws2.Cells.Item(62, 2).Value = 105.79
$ws2.Cells.Item(62, 3).Value = "forecast"
$ws2.Cells.Item(62, 4).Value = 46050.0
$ws2.Cells.Item(62, 5).Value = "ON"
$ws2.Cells.Item(62, 1).NumberFormat = $ws2.Cells.Item(61, 1).NumberFormat
$ws2.Cells.Item(62, 4).NumberFormat = $ws2.Cells.Item(61, 4).NumberFormat

# Row 63
$ws2.Cells.Item(63, 1).Value = 46050.27083333334
$ws2.Cells.Item(63, 2).Value = 97.3
$ws2.Cells.Item(63, 3).Value = "forecast"
$ws2.Cells.Item(63, 4).Value = 46050.0
$ws2.Cells.Item(63, 5).Value = "ON"
$ws2.Cells.Item(63, 1).NumberFormat = $ws2.Cells.Item(62, 1).NumberFormat
$ws2.Cells.Item(63, 4).NumberFormat = $ws2.Cells.Item(62, 4).NumberFormat

# Row 64
$ws2.Cells.Item(64, 1).Value = 46050.29166666666
$ws2.Cells.Item(64, 2).Value = 59.59837
$ws2.Cells.Item(64, 3).Value = "forecast"
$ws2.Cells.Item(64, 4).Value = 46050.0
$ws2.Cells.Item(64, 5).Value = "ON"
$ws2.Cells.Item(64, 1).NumberFormat = $ws2.Cells.Item(63, 1).NumberFormat
$ws2.Cells.Item(64, 4).NumberFormat = $ws2.Cells.Item(63, 4).NumberFormat

# Row 65
$ws2.Cells.Item(65, 1).Value = 46050.3125
$ws2.Cells.Item(65, 2).Value = 36.05919
$ws2.Cells.Item(65, 3).Value = "forecast"
$ws2.Cells.Item(65, 4).Value = 46050.0
$ws2.Cells.Item(65, 5).Value = "ON"
$ws2.Cells.Item(65, 1).NumberFormat = $ws2.Cells.Item(64, 1).NumberFormat
$ws2.Cells.Item(65, 4).NumberFormat = $ws2.Cells.Item(64, 4).NumberFormat

# Row 66
$ws2.Cells.Item(66, 1).Value = 46050.33333333334
$ws2.Cells.Item(66, 2).Value = 8.51369
$ws2.Cells.Item(66, 3).Value = "forecast"
$ws2.Cells.Item(66, 4).Value = 46050.0
$ws2.Cells.Item(66, 5).Value = "ON"
$ws2.Cells.Item(66, 1).NumberFormat = $ws2.Cells.Item(65, 1).NumberFormat
$ws2.Cells.Item(66, 4).NumberFormat = $ws2.Cells.Item(65, 4).NumberFormat

# Row 67
$ws2.Cells.Item(67, 1).Value = 46050.35416666666
$ws2.Cells.Item(67, 2).Value = 8.52165
$ws2.Cells.Item(67, 3).Value = "forecast"
$ws2.Cells.Item(67, 4).Value = 46050.0
$ws2.Cells.Item(67, 5).Value = "ON"
$ws2.Cells.Item(67, 1).NumberFormat = $ws2.Cells.Item(66, 1).NumberFormat
$ws2.Cells.Item(67, 4).NumberFormat = $ws2.Cells.Item(66, 4).NumberFormat

# Row 68
$ws2.Cells.Item(68, 1).Value = 46050.375
$ws2.Cells.Item(68, 2).Value = 0.72012
$ws2.Cells.Item(68, 3).Value = "forecast"
$ws2.Cells.Item(68, 4).Value = 46050.0
$ws2.Cells.Item(68, 5).Value = "ON"
$ws2.Cells.Item(68, 1).NumberFormat = $ws2.Cells.Item(67, 1).NumberFormat
$ws2.Cells.Item(68, 4).NumberFormat = $ws2.Cells.Item(67, 4).NumberFormat

# Row 69
$ws2.Cells.Item(69, 1).Value = 46050.39583333334
$ws2.Cells.Item(69, 2).Value = 0.51
$ws2.Cells.Item(69, 3).Value = "forecast"
$ws2.Cells.Item(69, 4).Value = 46050.0
$ws2.Cells.Item(69, 5).Value = "ON"
$ws2.Cells.Item(69, 1).NumberFormat = $ws2.Cells.Item(68, 1).NumberFormat
$ws2.Cells.Item(69, 4).NumberFormat = $ws2.Cells.Item(68, 4).NumberFormat

# Row 70
$ws2.Cells.Item(70, 1).Value = 46050.41666666666
$ws2.Cells.Item(70, 2).Value = -1.09497
$ws2.Cells.Item(70, 3).Value = "forecast"
$ws2.Cells.Item(70, 4).Value = 46050.0
$ws2.Cells.Item(70, 5).Value = "ON"
$ws2.Cells.Item(70, 1).NumberFormat = $ws2.Cells.Item(69, 1).NumberFormat
$ws2.Cells.Item(70, 4).NumberFormat = $ws2.Cells.Item(69, 4).NumberFormat

# Row 71
$ws2.Cells.Item(71, 1).Value = 46050.4375
$ws2.Cells.Item(71, 2).Value = -6.0
$ws2.Cells.Item(71, 3).Value = "forecast"
$ws2.Cells.Item(71, 4).Value = 46050.0
$ws2.Cells.Item(71, 5).Value = "ON"
$ws2.Cells.Item(71, 1).NumberFormat = $ws2.Cells.Item(70, 1).NumberFormat
$ws2.Cells.Item(71, 4).NumberFormat = $ws2.Cells.Item(70, 4).NumberFormat

# Row 72
$ws2.Cells.Item(72, 1).Value = 46050.45833333334
$ws2.Cells.Item(72, 2).Value = -6.0
$ws2.Cells.Item(72, 3).Value = "forecast"
$ws2.Cells.Item(72, 4).Value = 46050.0
$ws2.Cells.Item(72, 5).Value = "ON"
$ws2.Cells.Item(72, 1).NumberFormat = $ws2.Cells.Item(71, 1).NumberFormat
$ws2.Cells.Item(72, 4).NumberFormat = $ws2.Cells.Item(71, 4).NumberFormat

# Row 73
$ws2.Cells.Item(73, 1).Value = 46050.47916666666
$ws2.Cells.Item(73, 2).Value = -6.0
$ws2.Cells.Item(73, 3).Value = "forecast"
$ws2.Cells.Item(73, 4).Value = 46050.0
$ws2.Cells.Item(73, 5).Value = "ON"
$ws2.Cells.Item(73, 1).NumberFormat = $ws2.Cells.Item(72, 1).NumberFormat
$ws2.Cells.Item(73, 4).NumberFormat = $ws2.Cells.Item(72, 4).NumberFormat

# Row 74
$ws2.Cells.Item(74, 1).Value = 46050.5
$ws2.Cells.Item(74, 2).Value = -6.0
$ws2.Cells.Item(74, 3).Value = "forecast"
$ws2.Cells.Item(74, 4).Value = 46050.0
$ws2.Cells.Item(74, 5).Value = "ON"
$ws2.Cells.Item(74, 1).NumberFormat = $ws2.Cells.Item(73, 1).NumberFormat
$ws2.Cells.Item(74, 4).NumberFormat = $ws2.Cells.Item(73, 4).NumberFormat

# Row 75
$ws2.Cells.Item(75, 1).Value = 46050.52083333334
$ws2.Cells.Item(75, 2).Value = -5.89242
$ws2.Cells.Item(75, 3).Value = "forecast"
$ws2.Cells.Item(75, 4).Value = 46050.0
$ws2.Cells.Item(75, 5).Value = "ON"
$ws2.Cells.Item(75, 1).NumberFormat = $ws2.Cells.Item(74, 1).NumberFormat
$ws2.Cells.Item(75, 4).NumberFormat = $ws2.Cells.Item(74, 4).NumberFormat

# Row 76
$ws2.Cells.Item(76, 1).Value = 46050.54166666666
$ws2.Cells.Item(76, 2).Value = -1.14884
$ws2.Cells.Item(76, 3).Value = "forecast"
$ws2.Cells.Item(76, 4).Value = 46050.0
$ws2.Cells.Item(76, 5).Value = "ON"
$ws2.Cells.Item(76, 1).NumberFormat = $ws2.Cells.Item(75, 1).NumberFormat
$ws2.Cells.Item(76, 4).NumberFormat = $ws2.Cells.Item(75, 4).NumberFormat

# Row 77
$ws2.Cells.Item(77, 1).Value = 46050.5625
$ws2.Cells.Item(77, 2).Value = 0.00004
$ws2.Cells.Item(77, 3).Value = "forecast"
$ws2.Cells.Item(77, 4).Value = 46050.0
$ws2.Cells.Item(77, 5).Value = "ON"
$ws2.Cells.Item(77, 1).NumberFormat = $ws2.Cells.Item(76, 1).NumberFormat
$ws2.Cells.Item(77, 4).NumberFormat = $ws2.Cells.Item(76, 4).NumberFormat

# Row 78
$ws2.Cells.Item(78, 1).Value = 46050.58333333334
$ws2.Cells.Item(78, 2).Value = 0.0
$ws2.Cells.Item(78, 3).Value = "forecast"
$ws2.Cells.Item(78, 4).Value = 46050.0
$ws2.Cells.Item(78, 5).Value = "ON"
$ws2.Cells.Item(78, 1).NumberFormat = $ws2.Cells.Item(77, 1).NumberFormat
$ws2.Cells.Item(78, 4).NumberFormat = $ws2.Cells.Item(77, 4).NumberFormat

# Row 79
$ws2.Cells.Item(79, 1).Value = 46050.60416666666
$ws2.Cells.Item(79, 2).Value = -1.05659
$ws2.Cells.Item(79, 3).Value = "forecast"
$ws2.Cells.Item(79, 4).Value = 46050.0
$ws2.Cells.Item(79, 5).Value = "ON"
$ws2.Cells.Item(79, 1).NumberFormat = $ws2.Cells.Item(78, 1).NumberFormat
$ws2.Cells.Item(79, 4).NumberFormat = $ws2.Cells.Item(78, 4).NumberFormat

# Row 80
$ws2.Cells.Item(80, 1).Value = 46050.625
$ws2.Cells.Item(80, 2).Value = -5.94284
$ws2.Cells.Item(80, 3).Value = "forecast"
$ws2.Cells.Item(80, 4).Value = 46050.0
$ws2.Cells.Item(80, 5).Value = "ON"
$ws2.Cells.Item(80, 1).NumberFormat = $ws2.Cells.Item(79, 1).NumberFormat
$ws2.Cells.Item(80, 4).NumberFormat = $ws2.Cells.Item(79, 4).NumberFormat

# Row 81
$ws2.Cells.Item(81, 1).Value = 46050.64583333334
$ws2.Cells.Item(81, 2).Value = -5.94708
$ws2.Cells.Item(81, 3).Value = "forecast"
$ws2.Cells.Item(81, 4).Value = 46050.0
$ws2.Cells.Item(81, 5).Value = "ON"
$ws2.Cells.Item(81, 1).NumberFormat = $ws2.Cells.Item(80, 1).NumberFormat
$ws2.Cells.Item(81, 4).NumberFormat = $ws2.Cells.Item(80, 4).NumberFormat

# Row 82
$ws2.Cells.Item(82, 1).Value = 46050.66666666666
$ws2.Cells.Item(82, 2).Value = -4.88973
$ws2.Cells.Item(82, 3).Value = "forecast"
$ws2.Cells.Item(82, 4).Value = 46050.0
$ws2.Cells.Item(82, 5).Value = "OFF"
$ws2.Cells.Item(82, 1).NumberFormat = $ws2.Cells.Item(81, 1).NumberFormat
$ws2.Cells.Item(82, 4).NumberFormat = $ws2.Cells.Item(81, 4).NumberFormat

# Row 83
$ws2.Cells.Item(83, 1).Value = 46050.6875
$ws2.Cells.Item(83, 2).Value = -0.89805
$ws2.Cells.Item(83, 3).Value = "forecast"
$ws2.Cells.Item(83, 4).Value = 46050.0
$ws2.Cells.Item(83, 5).Value = "OFF"
$ws2.Cells.Item(83, 1).NumberFormat = $ws2.Cells.Item(82, 1).NumberFormat
$ws2.Cells.Item(83, 4).NumberFormat = $ws2.Cells.Item(82, 4).NumberFormat

# Row 84
$ws2.Cells.Item(84, 1).Value = 46050.70833333334
$ws2.Cells.Item(84, 2).Value = 9.67189
$ws2.Cells.Item(84, 3).Value = "forecast"
$ws2.Cells.Item(84, 4).Value = 46050.0
$ws2.Cells.Item(84, 5).Value = "OFF"
$ws2.Cells.Item(84, 1).NumberFormat = $ws2.Cells.Item(83, 1).NumberFormat
$ws2.Cells.Item(84, 4).NumberFormat = $ws2.Cells.Item(83, 4).NumberFormat

# Row 85
$ws2.Cells.Item(85, 1).Value = 46050.72916666666
$ws2.Cells.Item(85, 2).Value = 11.41325
$ws2.Cells.Item(85, 3).Value = "forecast"
$ws2.Cells.Item(85, 4).Value = 46050.0
$ws2.Cells.Item(85, 5).Value = "OFF"
$ws2.Cells.Item(85, 1).NumberFormat = $ws2.Cells.Item(84, 1).NumberFormat
$ws2.Cells.Item(85, 4).NumberFormat = $ws2.Cells.Item(84, 4).NumberFormat

# Row 86
$ws2.Cells.Item(86, 1).Value = 46050.75
$ws2.Cells.Item(86, 2).Value = 62.15067
$ws2.Cells.Item(86, 3).Value = "forecast"
$ws2.Cells.Item(86, 4).Value = 46050.0
$ws2.Cells.Item(86, 5).Value = "OFF"
$ws2.Cells.Item(86, 1).NumberFormat = $ws2.Cells.Item(85, 1).NumberFormat
$ws2.Cells.Item(86, 4).NumberFormat = $ws2.Cells.Item(85, 4).NumberFormat

# Row 87
$ws2.Cells.Item(87, 1).Value = 46050.77083333334
$ws2.Cells.Item(87, 2).Value = 68.77756
$ws2.Cells.Item(87, 3).Value = "forecast"
$ws2.Cells.Item(87, 4).Value = 46050.0
$ws2.Cells.Item(87, 5).Value = "OFF"
$ws2.Cells.Item(87, 1).NumberFormat = $ws2.Cells.Item(86, 1).NumberFormat
$ws2.Cells.Item(87, 4).NumberFormat = $ws2.Cells.Item(86, 4).NumberFormat

# Row 88
$ws2.Cells.Item(88, 1).Value = 46050.79166666666
$ws2.Cells.Item(88, 2).Value = 90.27302
$ws2.Cells.Item(88, 3).Value = "forecast"
$ws2.Cells.Item(88, 4).Value = 46050.0
$ws2.Cells.Item(88, 5).Value = "OFF"
$ws2.Cells.Item(88, 1).NumberFormat = $ws2.Cells.Item(87, 1).NumberFormat
$ws2.Cells.Item(88, 4).NumberFormat = $ws2.Cells.Item(87, 4).NumberFormat

# Row 89
$ws2.Cells.Item(89, 1).Value = 46050.8125
$ws2.Cells.Item(89, 2).Value = 91.62694
$ws2.Cells.Item(89, 3).Value = "forecast"
$ws2.Cells.Item(89, 4).Value = 46050.0
$ws2.Cells.Item(89, 5).Value = "OFF"
$ws2.Cells.Item(89, 1).NumberFormat = $ws2.Cells.Item(88, 1).NumberFormat
$ws2.Cells.Item(89, 4).NumberFormat = $ws2.Cells.Item(88, 4).NumberFormat

# Row 90
$ws2.Cells.Item(90, 1).Value = 46050.83333333334
$ws2.Cells.Item(90, 2).Value = 78.0
$ws2.Cells.Item(90, 3).Value = "forecast"
$ws2.Cells.Item(90, 4).Value = 46050.0
$ws2.Cells.Item(90, 5).Value = "OFF"
$ws2.Cells.Item(90, 1).NumberFormat = $ws2.Cells.Item(89, 1).NumberFormat
$ws2.Cells.Item(90, 4).NumberFormat = $ws2.Cells.Item(89, 4).NumberFormat

# Row 91
$ws2.Cells.Item(91, 1).Value = 46050.85416666666
$ws2.Cells.Item(91, 2).Value = 73.19
$ws2.Cells.Item(91, 3).Value = "forecast"
$ws2.Cells.Item(91, 4).Value = 46050.0
$ws2.Cells.Item(91, 5).Value = "OFF"
$ws2.Cells.Item(91, 1).NumberFormat = $ws2.Cells.Item(90, 1).NumberFormat
$ws2.Cells.Item(91, 4).NumberFormat = $ws2.Cells.Item(90, 4).NumberFormat

# Row 92
$ws2.Cells.Item(92, 1).Value = 46050.875
$ws2.Cells.Item(92, 2).Value = 64.89
$ws2.Cells.Item(92, 3).Value = "forecast"
$ws2.Cells.Item(92, 4).Value = 46050.0
$ws2.Cells.Item(92, 5).Value = "OFF"
$ws2.Cells.Item(92, 1).NumberFormat = $ws2.Cells.Item(91, 1).NumberFormat
$ws2.Cells.Item(92, 4).NumberFormat = $ws2.Cells.Item(91, 4).NumberFormat

# Row 93
$ws2.Cells.Item(93, 1).Value = 46050.89583333334
$ws2.Cells.Item(93, 2).Value = 64.89
$ws2.Cells.Item(93, 3).Value = "forecast"
$ws2.Cells.Item(93, 4).Value = 46050.0
$ws2.Cells.Item(93, 5).Value = "OFF"
$ws2.Cells.Item(93, 1).NumberFormat = $ws2.Cells.Item(92, 1).NumberFormat
$ws2.Cells.Item(93, 4).NumberFormat = $ws2.Cells.Item(92, 4).NumberFormat

# Row 94
$ws2.Cells.Item(94, 1).Value = 46050.91666666666
$ws2.Cells.Item(94, 2).Value = 57.09
$ws2.Cells.Item(94, 3).Value = "forecast"
$ws2.Cells.Item(94, 4).Value = 46050.0
$ws2.Cells.Item(94, 5).Value = "OFF"
$ws2.Cells.Item(94, 1).NumberFormat = $ws2.Cells.Item(93, 1).NumberFormat
$ws2.Cells.Item(94, 4).NumberFormat = $ws2.Cells.Item(93, 4).NumberFormat

# Row 95
$ws2.Cells.Item(95, 1).Value = 46050.9375
$ws2.Cells.Item(95, 2).Value = 57.09
$ws2.Cells.Item(95, 3).Value = "forecast"
$ws2.Cells.Item(95, 4).Value = 46050.0
$ws2.Cells.Item(95, 5).Value = "OFF"
$ws2.Cells.Item(95, 1).NumberFormat = $ws2.Cells.Item(94, 1).NumberFormat
$ws2.Cells.Item(95, 4).NumberFormat = $ws2.Cells.Item(94, 4).NumberFormat

# Row 96
$ws2.Cells.Item(96, 1).Value = 46050.95833333334
$ws2.Cells.Item(96, 2).Value = 57.06
$ws2.Cells.Item(96, 3).Value = "forecast"
$ws2.Cells.Item(96, 4).Value = 46050.0
$ws2.Cells.Item(96, 5).Value = "OFF"
$ws2.Cells.Item(96, 1).NumberFormat = $ws2.Cells.Item(95, 1).NumberFormat
$ws2.Cells.Item(96, 4).NumberFormat = $ws2.Cells.Item(95, 4).NumberFormat

# Row 97
$ws2.Cells.Item(97, 1).Value = 46050.97916666666
$ws2.Cells.Item(97, 2).Value = 57.06
$ws2.Cells.Item(97, 3).Value = "forecast"
$ws2.Cells.Item(97, 4).Value = 46050.0
$ws2.Cells.Item(97, 5).Value = "OFF"
$ws2.Cells.Item(97, 1).NumberFormat = $ws2.Cells.Item(96, 1).NumberFormat
$ws2.Cells.Item(97, 4).NumberFormat = $ws2.Cells.Item(96, 4).NumberFormat

